# Concurrent multi-task scraping results are merged into the "Hasil" sheet.
# Rows complete out of order (one goroutine per village/TPS group), which is
# reflected below by writing the rows in their original completion order.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hasil")
$ws.Activate()

# Row 3: UJONG MANGKI
$ws.Cells.Item(3,1).Value = "UJONG MANGKI"
$ws.Cells.Item(3,2).Value = 1101012002
$ws.Cells.Item(3,3).Value = 562
$ws.Cells.Item(3,4).Value = "[map[TPS1:[212 36 1]] map[TPS2:[235 26 7]]]"
$ws.Cells.Item(3,5).Value = 268
$ws.Cells.Item(3,6).Value = 294

# Row 5: GAMPONG DRIEN
$ws.Cells.Item(5,1).Value = "GAMPONG DRIEN"
$ws.Cells.Item(5,2).Value = 1101012004
$ws.Cells.Item(5,3).Value = 363
$ws.Cells.Item(5,4).Value = "[map[TPS1:[144 19 2]] map[TPS2:[153 14 4]]]"
$ws.Cells.Item(5,5).Value = 182
$ws.Cells.Item(5,6).Value = 181

# Row 6: DARUL IKHSAN
$ws.Cells.Item(6,1).Value = "DARUL IKHSAN"
$ws.Cells.Item(6,2).Value = 1101012015
$ws.Cells.Item(6,3).Value = 803
$ws.Cells.Item(6,4).Value = "[map[TPS1:[210 33 4]] map[TPS2:[207 31 1]] map[TPS3:[200 33 0]]]"
$ws.Cells.Item(6,5).Value = 271
$ws.Cells.Item(6,6).Value = 265
$ws.Cells.Item(6,7).Value = 267

# Row 7: PADANG BEURAHAN
$ws.Cells.Item(7,1).Value = "PADANG BEURAHAN"
$ws.Cells.Item(7,2).Value = 1101012016
$ws.Cells.Item(7,3).Value = 549
$ws.Cells.Item(7,4).Value = "[map[TPS1:[204 39 2]] map[TPS2:[203 36 1]]]"
$ws.Cells.Item(7,5).Value = 278
$ws.Cells.Item(7,6).Value = 271

# Row 8: GAMPONG BARO
$ws.Cells.Item(8,1).Value = "GAMPONG BARO"
$ws.Cells.Item(8,2).Value = 1101012017
$ws.Cells.Item(8,3).Value = 260
$ws.Cells.Item(8,4).Value = "[map[TPS1:[199 41 2]]]"
$ws.Cells.Item(8,5).Value = 260

# Row 9: FAJAR HARAPAN
$ws.Cells.Item(9,1).Value = "FAJAR HARAPAN"
$ws.Cells.Item(9,2).Value = 1101022001
$ws.Cells.Item(9,3).Value = 517
$ws.Cells.Item(9,4).Value = "[map[TPS1:[197 31 2]] map[TPS2:[196 26 1]]]"
$ws.Cells.Item(9,5).Value = 254
$ws.Cells.Item(9,6).Value = 263

# Row 10: KRUENG BATEE
$ws.Cells.Item(10,1).Value = "KRUENG BATEE"
$ws.Cells.Item(10,2).Value = 1101022002
$ws.Cells.Item(10,3).Value = 1102
$ws.Cells.Item(10,4).Value = "[map[TPS1:[212 16 2]] map[TPS2:[225 12 2]] map[TPS3:[229 10 1]] map[TPS4:[240 17 2]]]"
$ws.Cells.Item(10,5).Value = 271
$ws.Cells.Item(10,6).Value = 277
$ws.Cells.Item(10,7).Value = 266
$ws.Cells.Item(10,8).Value = 288

# Row 11: PASI KUALA ASAHAN
$ws.Cells.Item(11,1).Value = "PASI KUALA ASAHAN"
$ws.Cells.Item(11,2).Value = 1101022003
$ws.Cells.Item(11,3).Value = 549
$ws.Cells.Item(11,4).Value = "[map[TPS1:[235 15 2]] map[TPS2:[197 17 3]]]"
$ws.Cells.Item(11,5).Value = 282
$ws.Cells.Item(11,6).Value = 267

# Row 2: KEUDE BAKONGAN
$ws.Cells.Item(2,1).Value = "KEUDE BAKONGAN"
$ws.Cells.Item(2,2).Value = 1101012001
$ws.Cells.Item(2,3).Value = 873
$ws.Cells.Item(2,4).Value = "[map[TPS1:[186 44 8]] map[TPS2:[209 37 6]] map[TPS3:[202 38 7]]]"
$ws.Cells.Item(2,5).Value = 284
$ws.Cells.Item(2,6).Value = 296
$ws.Cells.Item(2,7).Value = 293

# Row 4: UJONG PADANG
$ws.Cells.Item(4,1).Value = "UJONG PADANG"
$ws.Cells.Item(4,2).Value = 1101012003
$ws.Cells.Item(4,3).Value = 476
$ws.Cells.Item(4,4).Value = "[map[TPS1:[176 27 3]] map[TPS2:[185 29 4]]]"
$ws.Cells.Item(4,5).Value = 230
$ws.Cells.Item(4,6).Value = 246

